# Add Ammar Rafaqat to the volunteer list (row 48).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 48 values ---
$ws.Range("A48").Value = "Ammar Rafaqat"
$ws.Range("B48").Value = "arrafaqat@gmail.com"
$ws.Range("C48").Value = 13
$ws.Range("D48").Value = "Acension Of Our Lord S.S."

# --- Formatting ---
# D48 uses the same "big contact-number" look already used elsewhere in the
# sheet (e.g. E46: Tahoma 15, dark-grey) -- copy that formatting over.
$ws.Range("E46").Copy()
$ws.Range("D48").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

# B48 (the email) gets a distinct teal Tahoma font. Start from a cell that
# is already Tahoma / teal (B45) so only the size needs to change, keeping
# the new style/font table addition minimal (one new font + one new xf).
$ws.Range("B45").Copy()
$ws.Range("B48").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false
$ws.Range("B48").Font.Size = 12

# Row 48 visually matches the taller rows (e.g. row 46) elsewhere in the
# sheet because of the larger font in D48.
$ws.Rows("48").RowHeight = 19
